{"js": "// 1) Insert a new \"Meta description\" paragraph right after the title\n//    (the first paragraph, styled Heading1).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst title = paragraphs.items[0];\n\n// Create a new empty paragraph right after the title and make sure it is\n// a plain body (\"Normal\") paragraph, not a copy of the Heading1 style.\nconst metaPara = title.insertParagraph(\"\", Word.InsertLocation.after);\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Insert the bold \"Meta description\" run, then the regular-text run that\n// follows it, as two separate runs (so only the label is bold).\nconst boldRun = metaPara.insertText(\"Meta description\", Word.InsertLocation.start);\nboldRun.font.bold = true;\nawait context.sync();\n\nconst restRun = metaPara.insertText(\n  \": Read our review of Anubis Wild Megaways, a unique online slot game with impressive graphics, bonus features, and a mysterious symbol. Play for free now.\",\n  Word.InsertLocation.end\n);\nrestRun.font.bold = false;\nawait context.sync();\n\n// 2) At the end of the document: remove the duplicated bold title\n//    paragraph, and rewrite the italic paragraph's text into the DALL-E\n//    image prompt (keeping its italic formatting).\nconst endParagraphs = body.paragraphs;\nendParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = endParagraphs.items;\nconst lastIndex = items.length - 1;\nconst italicPara = items[lastIndex];\nconst duplicateTitlePara = items[lastIndex - 1];\n\nduplicateTitlePara.delete();\nawait context.sync();\n\nitalicPara.insertText(\n  \"Prompt for DALLE: Create an engaging feature image for Anubis Wild Megaways in a cartoon style. The image should showcase a happy Maya warrior with glasses. Use vibrant and contrasting colors to make the image pop and attract the viewer's attention. Be creative with the design while incorporating symbols from the game, such as pyramids, sphinxes, and the Eye of Horus. The warrior should be holding a mobile device, indicating that the game is available to play online on different devices. The image should convey the excitement and adventure of playing Anubis Wild Megaways, inviting players to join the journey.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Insert a new \"Meta description\" paragraph right after the title\n#    (the first paragraph, styled Heading1).\n# ---------------------------------------------------------------------\n$title = $d.Paragraphs(1)\n$titleRange = $title.Range\n$titleRange.Collapse(0)            # wdCollapseEnd\n$titleRange.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaPara.Style = \"Normal\"\n\n$boldText = \"Meta description\"\n$boldRange = $metaPara.Range\n$boldStart = $boldRange.Start\n$boldRange.Text = $boldText\n$boldEnd = $boldStart + $boldText.Length\n\n# Format just the \"Meta description\" label as bold.\n$trueBoldRange = $d.Range($boldStart, $boldEnd)\n$trueBoldRange.Bold = 1\n\n# Append the (non-bold) rest of the meta description sentence.\n$restText = \": Read our review of Anubis Wild Megaways, a unique online slot game with impressive graphics, bonus features, and a mysterious symbol. Play for free now.\"\n$metaPara2 = $d.Paragraphs(2)\n$metaPara2.Range.InsertAfter($restText)\n\n$restRange = $d.Range($boldEnd, $boldEnd + $restText.Length)\n$restRange.Bold = 0\n\n# ---------------------------------------------------------------------\n# 2) At the end of the document: remove the duplicated bold title\n#    paragraph, and rewrite the italic paragraph's text into the DALL-E\n#    image prompt (keeping its italic formatting).\n# ---------------------------------------------------------------------\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs($count - 1)\n$duplicateTitlePara.Range.Delete()\n\n$count2 = $d.Paragraphs.Count\n$italicPara = $d.Paragraphs($count2)\n$italicRange = $italicPara.Range\n$italicTextRange = $d.Range($italicRange.Start, $italicRange.End - 1)\n$italicTextRange.Text = \"Prompt for DALLE: Create an engaging feature image for Anubis Wild Megaways in a cartoon style. The image should showcase a happy Maya warrior with glasses. Use vibrant and contrasting colors to make the image pop and attract the viewer's attention. Be creative with the design while incorporating symbols from the game, such as pyramids, sphinxes, and the Eye of Horus. The warrior should be holding a mobile device, indicating that the game is available to play online on different devices. The image should convey the excitement and adventure of playing Anubis Wild Megaways, inviting players to join the journey.\"\n"}
